$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A1:D4 table one column to the right (to B1:E4)
# by inserting a new blank column at column A.
$ws.Range("A1").EntireColumn.Insert()

# Update the selection to match the target state
$ws.Range("B10").Select()
